# Propellant Properties.xlsx - apply commit changes:
#  - rename Sheet1 -> Individual
#  - add Melting Point / Boiling Point rows to Individual
#  - add a new "Mixture" sheet with LOX/Ethanol mixture properties
#  - add a sheet-scoped defined name "tables" on Mixture!$D$1
#  - leave selection on Individual!E14 and Mixture!C14 (Mixture = active tab)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Individual sheet (was Sheet1): rename + append two new property rows
# ---------------------------------------------------------------------------
$individual = $wb.Worksheets.Item(1)
$individual.Name = "Individual"

$individual.Range("A13").Value = "Melting Point"
$individual.Range("B13").Value = -218.8
$individual.Range("C13").Value = -114.1

$individual.Range("A14").Value = "Boiling Point"
$individual.Range("B14").Value = -183
$individual.Range("C14").Value = 78.2

# ---------------------------------------------------------------------------
# 2. New "Mixture" sheet, placed right after Individual
# ---------------------------------------------------------------------------
$mixture = $wb.Worksheets.Add($null, $individual)
$mixture.Name = "Mixture"

# Header row (bold-looking italic/underline, size 16 - same look as the
# "Properties/Fuel/Oxidizer" header row on Individual)
$mixture.Range("B1").Value = "Oxidizer - LOX"
$mixture.Range("C1").Value = "Fuel - Ethanol"
$headerRange = $mixture.Range("B1:D1")
$headerRange.Font.Italic = $true
$headerRange.Font.Underline = $true
$headerRange.Font.Size = 16
$mixture.Rows.Item(1).RowHeight = 21

# Data rows
$mixture.Range("A2").Value = "Hypergolic"
$mixture.Range("B2").Value = "No"

$mixture.Range("A3").Value = "Mixture Ratio"
$mixture.Range("B3").Value = 1.29

$mixture.Range("A4").Value = "Specific Impulse (Sea Level)"
$mixture.Range("B4").Value = 269

$mixture.Range("A5").Value = "Density Impulse"
$mixture.Range("B5").Value = 264

# Merge + center the value cells (B:C) for each data row
foreach ($r in 2..5) {
    $rowRange = $mixture.Range("B" + $r + ":C" + $r)
    $rowRange.Merge() | Out-Null
    $rowRange.HorizontalAlignment = -4108   # xlCenter
}

# Column widths (bestFit-like, in characters)
$mixture.Columns.Item(1).ColumnWidth = 23.498697916666668
$mixture.Columns.Item(2).ColumnWidth = 16.498697916666668
$mixture.Columns.Item(3).ColumnWidth = 16.330729166666668

# ---------------------------------------------------------------------------
# 3. Sheet-scoped defined name "tables" on Mixture!$D$1
# ---------------------------------------------------------------------------
$mixture.Names.Add("tables", "=Mixture!`$D`$1")

# ---------------------------------------------------------------------------
# 4. Final selection / active sheet state
# ---------------------------------------------------------------------------
$individual.Range("E14").Select() | Out-Null
$mixture.Range("C14").Select() | Out-Null
